$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old lone data row (row 4); its data is superseded by the new table below
[void]$ws.Range("A4:E4").ClearContents()

# New column header
$ws.Range("F3").Value = "Epochs_needed"

# New statistics table (rows 5-11) for the various OneHot scaling factors
$data = @(
  @(0.5,  0.000205, 0.000543,  122.332,    0.06103,  13900),
  @(1.5,  0.000101, 0.0008,    156.1654,   0.6389,   10200),
  @(3,    0.000259, 0.0003821, 121.8057,   0.61251,  11600),
  @(5,    0.000367, 0.00188,   268.127,    0.712853, 33900),
  @(8,    0.0007,   0.000157,  33.3766,    0.6698,   30800),
  @(10,   0.001558, 0.001754,  254.02193,  0.732331, 11100),
  @(15,   0.000916, 0.0002433, 39.4652,    0.6543,   27300)
)

$r = 5
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Highlight/format the standout B9 cell (custom font + left/center alignment)
$ws.Range("B9").Font.Color = 0
$ws.Range("B9").Font.Name = "Calibri "
$ws.Range("B9").HorizontalAlignment = -4131
$ws.Range("B9").VerticalAlignment = -4108

# Move the active selection like the author left it
[void]$ws.Range("F12").Select()

# Page setup tweak (A4 portrait) from the xls updates
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "applied onehot statistics update"
